$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 in the original workbook only holds a stray "mapping-ano.xlsx" value in H5
# and is not part of the real metadata table. Drop it before inserting the new
# "slug" row so it doesn't get dragged down to row 6.
$ws.Rows.Item(5).Delete()

# Insert a brand-new row 2 that will carry machine-friendly "slug" column
# identifiers; this pushes the old rows 2-4 (iaest-measure/dim, medida/dim,
# xsd:double/URI-*) down to rows 3-5, matching the target layout.
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "superficie-agricola-utilizada-sau"
$ws.Range("B2").Value = "ccaa-nombre"
$ws.Range("C2").Value = "ccaa-codigo"
$ws.Range("D2").Value = "provincia-codigo"
$ws.Range("E2").Value = "superficie-regada-sobre-sau"
$ws.Range("F2").Value = "municipio-codigo"
$ws.Range("G2").Value = "provincia-nombre"
$ws.Range("H2").Value = "ano"
$ws.Range("I2").Value = "superficie-regada-has"
$ws.Range("J2").Value = "municipio-nombre"
